# PlayerPerformance_4704.xlsx edit:
#  1. Insert a new "Player Info" worksheet before the existing "ODI Batting"
#     sheet and populate it with the player's basic info.
#  2. On "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#     full howstat.com scorecard URLs with the bare match-code number.

$wb = $excel.ActiveWorkbook

# --- existing sheet reference (grab it before inserting, by name) ---------
$batting = $wb.Worksheets.Item("ODI Batting")

# --- 1. new "Player Info" sheet, inserted before "ODI Batting" ------------
$info = $wb.Worksheets.Add($batting)
$info.Name = "Player Info"

# NOTE: inserting a sheet "before" $batting shifts what that positional
# reference resolves to (it now points at the newly inserted sheet), so
# re-resolve the "ODI Batting" sheet by name before touching it again.
$batting = $wb.Worksheets.Item("ODI Batting")

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

$headerRange = $info.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Player id must be stored as text (matches the source sheet convention of
# keeping all data as text), so force a text number format before writing,
# then drop back to the default style so no stray format is attached.
$idCell = $info.Range("A2")
$idCell.NumberFormat = "@"
$idCell.Value = "4704"
$idCell.Style = "Normal"

$info.Range("B2").Value = "Wedagedara Sadeera Rashen Samarawickrama"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Does Not Bowl | Unknown"

# --- 2. "ODI Batting" sheet edits ------------------------------------------
$batting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{
    2 = "4084"
    3 = "4087"
    4 = "4099"
    5 = "4211"
    6 = "4212"
    7 = "4215"
    8 = "4375"
}

foreach ($row in $matchCodes.Keys) {
    $cell = $batting.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$row]
    $cell.Style = "Normal"
}

Write-Output "Player Info sheet added; ODI Batting MATCH_CODE column updated."
